$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Shift the 90-day rolling window up by one row: row r takes what used to be
# in row r+1 (dropping the oldest day, row 2, and leaving row 91 to be filled
# in with the newest day below). Column A (date, stored as text) is moved via
# Copy so the shared-string/text typing is preserved instead of being
# reinterpreted as a date serial; column C (page count) is a plain number.
for ($r = 2; $r -le 90; $r++) {
    $ws.Cells.Item($r + 1, 1).Copy($ws.Cells.Item($r, 1))
    $ws.Cells.Item($r + 1, 3).Copy($ws.Cells.Item($r, 3))
}

# New newest day appended at the end of the window.
$lastRow = 91
$dateCell = $ws.Cells.Item($lastRow, 1)
$dateCell.Formula = '="2026-01-20"'
$dateCell.Copy()
$dateCell.PasteSpecial(-4163) # xlPasteValues - keep as literal text, not a formula

$ws.Cells.Item($lastRow, 3).Value = 26.0
